$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 9 (shifts old rows 9-12 down to 13-16)
$ws.Rows("9:12").Insert()

# The two "template" rows (current rows 7 and 8, siniestro group 301) get duplicated
# into the newly inserted rows 9-12 (9=copy of 7, 10=copy of 8, 11=copy of 7, 12=copy of 8)
$ws.Range("A7:AH8").Copy($ws.Range("A9:AH10"))
$ws.Range("A7:AH8").Copy($ws.Range("A11:AH12"))

# Update the siniestro id in rows 7 and 8 (and the freshly pasted copies) to 301
$ws.Range("A7").Value = 301
$ws.Range("A8").Value = 301
$ws.Range("A9").Value = 301
$ws.Range("A10").Value = 301
$ws.Range("A11").Value = 301
$ws.Range("A12").Value = 301

# Restore the selection Excel leaves active after this kind of edit
$ws.Range("A7").Select() | Out-Null
